$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.313.81"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.682.39"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.43"
$ws.Range("E5").Value = "  +3.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.31"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.682.72"
$ws.Range("E9").Value = "  +4.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.106"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.152"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.34"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.166.64"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.219.23"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000144"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.693.10"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.43"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.19"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.39"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.87"
$ws.Range("E22").Value = "  +4.06%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.59"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  +4.01%  "
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.53"
$ws.Range("E28").Value = "  +4.11%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "539.38"
$ws.Range("E30").Value = "  +19.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.86"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.02"
$ws.Range("E32").Value = "  +5.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("E33").Value = "  +10.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0810"
$ws.Range("E34").Value = "  +2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "173.63"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("E36").Value = "  +15.12%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.406"
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.23"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.82"
$ws.Range("E40").Value = "  +7.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.44"
$ws.Range("E41").Value = "  +13.14%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.10"
$ws.Range("E44").Value = "  +5.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0563"
$ws.Range("E45").Value = "  +5.68%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0961"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0239"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.01"
$ws.Range("E49").Value = "  +6.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.32"
$ws.Range("E51").Value = "  -0.73%  "
